# Adds the "Progreso Día 4" block (rows 34-37) to Hoja1, mirroring the
# existing "Progreso Día 3" block's layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from the analogous rows of the "Progreso Día 3" block
# so the new rows pick up the same cell styles (title / header / data).
$ws.Range("A25:E25").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)

$ws.Range("A26:E26").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)

$ws.Range("A27:E27").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)

$ws.Range("A27:E27").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)

# --- Row 34: section title (merged A34:E34)
$ws.Range("A34").Value = "Progreso Día 4"
$ws.Range("A34:E34").Merge()

# --- Row 35: column headers
$ws.Range("A35").Value = "Fecha"
$ws.Range("B35").Value = "Tarea"
$ws.Range("C35").Value = "Descripción del Progreso"
$ws.Range("D35").Value = "Archivos Modificados"
$ws.Range("E35").Value = "Observaciones"

# --- Row 36: "Comentarios Detallados" entry
$ws.Range("A36").Value = "2024-07-22"
$ws.Range("B36").Value = "Comentarios Detallados"

$c36Text = "Se añadieron comentarios detallados a los archivos de código EditCountry.jsx y NotFound.jsx, explicando cada línea del código para mejorar la comprensión."
$ws.Range("C36").Value = $c36Text

$part1 = "Se añadieron comentarios detallados a los archivos de código "
$part2 = "EditCountry.jsx"
$part3 = " y "
$part4 = "NotFound.jsx"
$part5 = ", explicando cada línea del código para mejorar la comprensión."

$pos2 = $part1.Length + 1
$pos3 = $pos2 + $part2.Length
$pos4 = $pos3 + $part3.Length
$pos5 = $pos4 + $part4.Length

$run2 = $ws.Range("C36").Characters($pos2, $part2.Length)
$run2.Font.Name = "Arial Unicode MS"
$run2.Font.Size = 10

$run3 = $ws.Range("C36").Characters($pos3, $part3.Length)
$run3.Font.Name = "Aptos Narrow"
$run3.Font.Size = 11

$run4 = $ws.Range("C36").Characters($pos4, $part4.Length)
$run4.Font.Name = "Arial Unicode MS"
$run4.Font.Size = 10

$run5 = $ws.Range("C36").Characters($pos5, $part5.Length)
$run5.Font.Name = "Aptos Narrow"
$run5.Font.Size = 11

$ws.Range("D36").Value = "Todos los JS y JSX"
$ws.Range("E36").Value = "Comentarios detallados añadidos para una mejor comprensión del código."

# --- Row 37: "Documentación Word" entry
$ws.Range("A37").Value = "2024-07-22"
$ws.Range("B37").Value = "Documentación Word"
$ws.Range("C37").Value = "Se creó un documento Word con explicaciones detalladas del código, incluyendo fragmentos específicos para ilustrar cada punto."
$ws.Range("D37").Value = "N/A"
$ws.Range("E37").Value = "Documento creado para proporcionar una comprensión clara y completa del funcionamiento del código."

# --- Row heights to match the authored layout
$ws.Rows.Item(34).RowHeight = 42
$ws.Rows.Item(35).RowHeight = 31.5
$ws.Rows.Item(36).RowHeight = 90
$ws.Rows.Item(37).RowHeight = 135

# --- Restore the view/selection state shown after the edit
$ws.Range("H37").Select()
